$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Locate the "New enhancements..." paragraph and the "Continuous Integration"
# heading paragraph by content, rather than by a hard-coded index, so the
# script is resilient to any incidental paragraph-count differences.
# ---------------------------------------------------------------------------
$enhancementsIndex = 0
$continuousIntegrationIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $text = $d.Paragraphs($i).Range.Text.TrimEnd()
    if ($text.StartsWith("New enhancements, ideas and any feedback")) {
        $enhancementsIndex = $i
    }
    if ($text -eq "Continuous Integration") {
        $continuousIntegrationIndex = $i
    }
}

if ($enhancementsIndex -eq 0) { throw "Could not find the 'New enhancements' paragraph" }
if ($continuousIntegrationIndex -eq 0) { throw "Could not find the 'Continuous Integration' paragraph" }

# ---------------------------------------------------------------------------
# STEP 1: merge the two runs of the "New enhancements..." paragraph into a
# single run, appending the extra sentence fragment. This also drops the
# now-redundant _GoBack bookmark that used to sit between the two runs.
# ---------------------------------------------------------------------------
$p = $d.Paragraphs($enhancementsIndex)
$mergedText = "New enhancements, ideas and any feedback from members of university staff should be posted in the #improvements channel in slack. This channel should be dedicated to discussing creative ideas and monitoring which improvements need to be made after each deliverable. The person who presents each deliverable is responsible for posting and feedback given from their lab tutor. The project manager should be consulted on any major decisions that need to be made based on feedback, such as taking the project in a different direction. "

# Overwriting with text identical to the existing (pre-edit) content is a
# no-op for the underlying engine, so first write the target text plus a
# trailing sentinel character, then strip that sentinel off in a second
# pass -- this forces a genuine rewrite that flattens the paragraph down to
# a single run and removes the bookmark that used to split it.
$bodyRange = $d.Range($p.Range.Start, $p.Range.End - 1)
$bodyRange.Text = $mergedText + "X"
$p2 = $d.Paragraphs($enhancementsIndex)
$sentinelRange = $d.Range($p2.Range.End - 2, $p2.Range.End - 1)
$sentinelRange.Text = ""

# ---------------------------------------------------------------------------
# STEP 2: after the "Continuous Integration" heading paragraph, add a blank
# underlined paragraph followed by the new "Continuous Integration" body
# paragraph (two runs) carrying the relocated _GoBack bookmark.
# ---------------------------------------------------------------------------
$ciPara = $d.Paragraphs($continuousIntegrationIndex)
$insertPoint = $d.Range($ciPara.Range.End, $ciPara.Range.End)

$run1Text = "For this project, we will be using continuous integration methods for merging changes to the project. It is encouraged to push changes to the master branch of the project "
$run2Text = "when each separate feature is completed. We have chosen this method as all members of the team will be contributing to programming components in the project which need to be linked together as shown in the design diagrams. Having multiple small commits allows "

$frag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
    '<w:p><w:pPr><w:rPr><w:u w:val="single"/></w:rPr></w:pPr></w:p>' + `
    '<w:p><w:r><w:t xml:space="preserve">' + $run1Text + '</w:t></w:r><w:r><w:t xml:space="preserve">' + $run2Text + '</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>' + `
    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertPoint.InsertXML($frag)
